$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The monthly "Periodo Mora" (column E) rows are reordered to descending
# (most recent period first) and the "Valor Mora" (column F) values follow
# the same row-for-row swap. Net effect vs the original layout: row 16 and
# row 22 trade places (2404/52000 <-> 2410/15600), row 17 and row 21 trade
# places (2405/52000 <-> 2409/52000), row 18 and row 20 trade places
# (2406/52000 <-> 2408/52000), and row 19 (2407/52000) is unchanged.

$periodos = @("2410", "2409", "2408", "2407", "2406", "2405", "2404")
$valores  = @(15600, 52000, 52000, 52000, 52000, 52000, 52000)

for ($i = 0; $i -lt 7; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
    $ws.Range("F$row").Value = $valores[$i]
}
